# Natmi following Dr Hou advice
# Rebuild the Gdf9 -> Acvr2a sending/target cluster matrix (ECs, FAPs, sCs)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Gdf9", "Acvr2a", "ECs", 3, 1, 0.9029553333333333, 2.708866, 0.2353208431348719, 0.2353208431348719, 3, 1, 14.61878266666667, 43.856348, 0.2662829816142094, 0.2662829816142094, 13.20010777570756, 118.800969981368, 0.06266193574592337, 0.06266193574592337),
    @("ECs", "Gdf9", "Acvr2a", "FAPs", 3, 1, 0.9029553333333333, 2.708866, 0.2353208431348719, 0.2353208431348719, 3, 1, 27.084169, 81.25250700000001, 0.4933415757187404, 0.4933415757187404, 24.45579484745134, 220.102153627062, 0.1160935555516202, 0.1160935555516203),
    @("ECs", "Gdf9", "Acvr2a", "sCs", 3, 1, 0.9029553333333333, 2.708866, 0.2353208431348719, 0.2353208431348719, 3, 1, 13.19647366666667, 39.589421, 0.2403754426670501, 0.2403754426670501, 11.91582627850956, 107.242436506586, 0.0565653518373283, 0.0565653518373283),
    @("FAPs", "Gdf9", "Acvr2a", "ECs", 3, 1, 1.446898, 4.340694, 0.3770787376970584, 0.3770787376970584, 3, 1, 14.61878266666667, 43.856348, 0.2662829816142094, 0.2662829816142094, 21.15188740283467, 190.366986625512, 0.1004096505772951, 0.1004096505772951),
    @("FAPs", "Gdf9", "Acvr2a", "FAPs", 3, 1, 1.446898, 4.340694, 0.3770787376970584, 0.3770787376970584, 3, 1, 27.084169, 81.25250700000001, 0.4933415757187404, 0.4933415757187404, 39.18802995776201, 352.6922696198581, 0.1860286186255004, 0.1860286186255004),
    @("FAPs", "Gdf9", "Acvr2a", "sCs", 3, 1, 1.446898, 4.340694, 0.3770787376970584, 0.3770787376970584, 3, 1, 13.19647366666667, 39.589421, 0.2403754426670501, 0.2403754426670501, 19.09395135535267, 171.845562198174, 0.0906404684942629, 0.09064046849426288),
    @("sCs", "Gdf9", "Acvr2a", "ECs", 2, 0.6666666666666666, 1.487271, 4.461812999999999, 0.3876004191680696, 0.3876004191680697, 3, 1, 14.61878266666667, 43.856348, 0.2662829816142094, 0.2662829816142094, 21.74209151543599, 195.678823638924, 0.103211395290991, 0.103211395290991),
    @("sCs", "Gdf9", "Acvr2a", "FAPs", 2, 0.6666666666666666, 1.487271, 4.461812999999999, 0.3876004191680696, 0.3876004191680697, 3, 1, 27.084169, 81.25250700000001, 0.4933415757187404, 0.4933415757187404, 40.281499112799, 362.533492015191, 0.1912194015416197, 0.1912194015416198),
    @("sCs", "Gdf9", "Acvr2a", "sCs", 2, 0.6666666666666666, 1.487271, 4.461812999999999, 0.3876004191680696, 0.3876004191680697, 3, 1, 13.19647366666667, 39.589421, 0.2403754426670501, 0.2403754426670501, 19.626732586697, 176.640593280273, 0.09316962233545892, 0.09316962233545893)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $rows[$i]
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
